$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 209, shifting existing rows 209-220 down to 210-221
$ws.Rows.Item(209).Insert()

# Populate the newly inserted row 209 with the new data record
$ws.Range("A209").Value = 3
$ws.Range("B209").Value = "Femacal de La Calera"
$ws.Range("C209").Value = "Coquimbo"
$ws.Range("D209").Value = 44931
$ws.Range("E209").Value = 5
$ws.Range("F209").Value = 100112030
$ws.Range("G209").Value = "Poroto granado"
$ws.Range("H209").Value = "Sin especificar"
$ws.Range("I209").Value = "Primera"
$ws.Range("J209").Value = 78
$ws.Range("K209").Value = 44000
$ws.Range("L209").Value = 45000
$ws.Range("M209").Value = 44513
$ws.Range("N209").Value = "$/malla 25 kilos"
$ws.Range("O209").Value = "Provincia de Quillota"
$ws.Range("P209").Value = 1781
$ws.Range("Q209").Value = 25
$ws.Range("R209").Value = "Hortaliza"
